$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 6, column A) onto the
# new row's column A cell so the new participant number picks up the same
# bold/bordered/centered style used by the other "N_Ano" cells.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Fill in the new participant's answers (row 7)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "15h26"
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = "Autre"
$ws.Range("E7").Value = "4h"
$ws.Range("F7").Value = "Non"
$ws.Range("G7").Value = "Assez stressé(e)"
$ws.Range("H7").Value = "Non"
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "Non"
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = "Compétent"
$ws.Range("M7").Value = "2 ans"
$ws.Range("N7").Value = "Légèrement passionné(e)"
$ws.Range("O7").Value = "Très bruyant"
$ws.Range("P7").Value = "test55555"
